# Apply updated Betfair back/lay odds values per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.48
$ws.Range("Q2").Value = 1.63
# Row 3
$ws.Range("S3").Value = 2.66
# Row 4
$ws.Range("F4").Value = 1.53
$ws.Range("H4").Value = 6.8
$ws.Range("P4").Value = 2.08
$ws.Range("S4").Value = 3.1
$ws.Range("U4").Value = 1.95
$ws.Range("V4").Value = 1.15
$ws.Range("X4").Value = 19.5
$ws.Range("AO4").Value = 160
# Row 5
$ws.Range("F5").Value = 4.6
$ws.Range("G5").Value = 5.1
$ws.Range("H5").Value = 1.85
$ws.Range("I5").Value = 1.91
$ws.Range("N5").Value = 3.2
$ws.Range("P5").Value = 1.76
$ws.Range("R5").Value = 1.29
$ws.Range("T5").Value = 2.02
$ws.Range("U5").Value = 1.86
$ws.Range("V5").Value = 2.1
$ws.Range("W5").Value = 1.25
# Row 6
$ws.Range("F6").Value = 1.42
$ws.Range("H6").Value = 9.4
$ws.Range("I6").Value = 12
$ws.Range("K6").Value = 5.3
$ws.Range("W6").Value = 2.96
# Row 7
$ws.Range("F7").Value = 1.96
$ws.Range("J7").Value = 3.55
$ws.Range("K7").Value = 3.65
$ws.Range("Q7").Value = 2.18
$ws.Range("T7").Value = 1.94
# Row 8
$ws.Range("J8").Value = 5.8
$ws.Range("P8").Value = 2.86
$ws.Range("U8").Value = 2.24
$ws.Range("AO8").Value = 80
# Row 9
$ws.Range("G9").Value = 2.44
$ws.Range("H9").Value = 3.15
$ws.Range("I9").Value = 3.65
$ws.Range("M9").Value = 1.06
$ws.Range("Q9").Value = 1.82
$ws.Range("V9").Value = 1.38
$ws.Range("W9").Value = 1.69
# Row 10
$ws.Range("N10").Value = 2.96
$ws.Range("Q10").Value = 2.54
$ws.Range("T10").Value = 2.04
$ws.Range("U10").Value = 1.92
$ws.Range("W10").Value = 1.58
# Row 11
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 3.2
$ws.Range("Q11").Value = 2.16
# Row 12
$ws.Range("J12").Value = 3.25
$ws.Range("L12").Value = 1.45
$ws.Range("S12").Value = 3.95
# Row 13
$ws.Range("F13").Value = 2.38
$ws.Range("G13").Value = 2.54
$ws.Range("I13").Value = 3.6
$ws.Range("J13").Value = 3.2
$ws.Range("K13").Value = 3.45
$ws.Range("V13").Value = 1.4
$ws.Range("W13").Value = 1.65
# Row 16
$ws.Range("S16").Value = 1.84
# Row 17
$ws.Range("G17").Value = 4.3
$ws.Range("J17").Value = 4
$ws.Range("L17").Value = 1.27
$ws.Range("N17").Value = 5.6
$ws.Range("P17").Value = 2.58
$ws.Range("Q17").Value = 1.52
$ws.Range("R17").Value = 1.65
$ws.Range("S17").Value = 2.28
$ws.Range("T17").Value = 1.55
$ws.Range("U17").Value = 2.54
$ws.Range("W17").Value = 1.3
$ws.Range("AH17").Value = 18
$ws.Range("AM17").Value = 60
$ws.Range("AN17").Value = 32
# Row 18
$ws.Range("F18").Value = 1.83
$ws.Range("G18").Value = 2.08
$ws.Range("H18").Value = 4.6
$ws.Range("I18").Value = 5.7
$ws.Range("J18").Value = 3.25
$ws.Range("K18").Value = 3.85
$ws.Range("L18").Value = 1.42
$ws.Range("T18").Value = 1.95
$ws.Range("U18").Value = 1.84
$ws.Range("V18").Value = 1.22
$ws.Range("W18").Value = 1.93
$ws.Range("X18").Value = 13
$ws.Range("Y18").Value = 970
$ws.Range("Z18").Value = 44
$ws.Range("AA18").Value = 150
$ws.Range("AD18").Value = 23
$ws.Range("AE18").Value = 90
$ws.Range("AF18").Value = 13
$ws.Range("AI18").Value = 100
$ws.Range("AJ18").Value = 26
$ws.Range("AM18").Value = 170
# Row 19
$ws.Range("F19").Value = 2.44
$ws.Range("G19").Value = 2.56
$ws.Range("H19").Value = 3.45
$ws.Range("I19").Value = 3.7
$ws.Range("N19").Value = 2.64
$ws.Range("Q19").Value = 2.56
$ws.Range("V19").Value = 1.37
$ws.Range("W19").Value = 1.64
$ws.Range("AC19").Value = 6.8
$ws.Range("AN19").Value = 38
$ws.Range("AO19").Value = 95
# Row 20
$ws.Range("F20").Value = 2.68
$ws.Range("I20").Value = 2.94
$ws.Range("Q20").Value = 2.16
# Row 21
$ws.Range("G21").Value = 2.94
$ws.Range("J21").Value = 3.1
$ws.Range("K21").Value = 3.15
$ws.Range("O21").Value = 1.48
$ws.Range("P21").Value = 1.65
$ws.Range("Q21").Value = 2.48
$ws.Range("Z21").Value = 18.5
$ws.Range("AF21").Value = 18
$ws.Range("AH21").Value = 20
$ws.Range("AM21").Value = 140
# Row 22
$ws.Range("S22").Value = 2.32
$ws.Range("U22").Value = 2.76
$ws.Range("AG22").Value = 10.5
$ws.Range("AH22").Value = 14
# Row 23
$ws.Range("H23").Value = 3.55
$ws.Range("J23").Value = 3.45
$ws.Range("P23").Value = 1.8
$ws.Range("R23").Value = 1.29
$ws.Range("W23").Value = 1.75
$ws.Range("AB23").Value = 8.800000000000001
# Row 24
$ws.Range("H24").Value = 2.62
$ws.Range("Q24").Value = 1.98
$ws.Range("T24").Value = 1.73
$ws.Range("W24").Value = 1.51
$ws.Range("AC24").Value = 7.6
